$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing Description text for several rows ---
$ws.Range("D17").Value = "2016 January Lobbyist Report"
$ws.Range("D18").Value = "2016 May Lobbyist Report"
$ws.Range("D20").Value = "2016 October Lobbyist Report"
$ws.Range("D22").Value = "2017 January Lobbyist Report"
$ws.Range("D23").Value = "2017 January 48 Hour Report"
$ws.Range("D24").Value = "2017 May Lobbyist Report"
$ws.Range("D26").Value = "2017 October Lobbyist Report"

# --- Append three new rows (27-29) of 2018 filing data ---
# Copy formatting from the last existing data row (26) down into the new rows
$ws.Range("A26:G26").Copy()
$ws.Range("A27:G29").PasteSpecial(-4122)

# Row 27
$ws.Cells.Item(27, 1).Value = 26
$ws.Cells.Item(27, 2).Value = 43115.4993055556
$ws.Cells.Item(27, 3).Value = 43100.4993055556
$ws.Cells.Item(27, 4).Value = "2018 January Lobbyist Report"
$ws.Cells.Item(27, 5).Value = 1
$ws.Cells.Item(27, 6).Value = 1
$ws.Cells.Item(27, 7).Value = 43011

# Row 28
$ws.Cells.Item(28, 1).Value = 27
$ws.Cells.Item(28, 2).Value = 43148
$ws.Cells.Item(28, 3).Value = 43146
$ws.Cells.Item(28, 4).Value = "2018 January 48 Hour Report"
$ws.Cells.Item(28, 5).Value = 2
$ws.Cells.Item(28, 6).Value = 1
$ws.Cells.Item(28, 7).Value = 43116

# Row 29
$ws.Cells.Item(29, 1).Value = 28
$ws.Cells.Item(29, 2).Value = 43229.9993055556
$ws.Cells.Item(29, 3).Value = 43227.9993055556
$ws.Cells.Item(29, 4).Value = "2018 May Lobbyist Report"
$ws.Cells.Item(29, 5).Value = 1
$ws.Cells.Item(29, 6).Value = 1
$ws.Cells.Item(29, 7).Value = 43101
